$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 11.244096314590454
$ws.Range("C2").Value = 15.395064886753421
$ws.Range("D2").Value = 9.4304240103828096
$ws.Range("E2").Value = 14.524984440938098

$ws.Range("B3").Value = 14.207823722985925
$ws.Range("C3").Value = 16.513892442244266
$ws.Range("D3").Value = 16.832541100581402
$ws.Range("E3").Value = 15.615825002300481

$ws.Activate()
$ws.Range("B1:E3").Select()
